$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.723.17"
$ws.Range("E2").Value = "  +2.00%  "
$ws.Range("D3").Value = "1.796.66"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.54"
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.555"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.44"
$ws.Range("E8").Value = "  +3.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.287"
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("E10").Value = "  +9.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0933"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "2.053.30"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.825.75"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.15"
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.640"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "34.722.30"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.34"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "254.31"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "0.0₃0812"
$ws.Range("E20").Value = "  +8.98%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.79"
$ws.Range("E22").Value = "  +3.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.24"
$ws.Range("E23").Value = "  -1.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.14"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.63"
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.44"
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.14"
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0536"
$ws.Range("E30").Value = "  +3.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.81"
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.62"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("D35").Value = "1.441.56"
$ws.Range("E35").Value = "  -2.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.07"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("E37").Value = "  +2.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.636"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "85.06"
$ws.Range("E39").Value = "  +1.24%  "
$ws.Range("B40").Value = "Swop.fi"
$ws.Range("C40").Value = "https://coinranking.com/coin/yrCr2HW2c+swopfi-swop"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "372.93"
$ws.Range("E40").Value = "  +606.54%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.79"
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.935"
$ws.Range("E42").Value = "  +3.59%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.33"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.13"
$ws.Range("E44").Value = "  +3.38%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.00"
$ws.Range("E45").Value = "  +5.30%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.06"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0495"
$ws.Range("E47").Value = "  -3.94%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.950.07"
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.32"
$ws.Range("E49").Value = "  +8.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "12.03"
$ws.Range("E50").Value = "  +1.63%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -0.08%  "
